$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "ECs" sending-cluster row (old row 2).
# FAPs/MuSCs rows shift up to become rows 2 and 3.
$ws.Rows.Item(2).Delete()

# Update TPM-derived values for the FAPs -> Thpo/Mpl -> ECs edge (now row 2)
$ws.Range("H2").Value = 7.606659
$ws.Range("I2").Value = 0.8548859591470492
$ws.Range("J2").Value = 0.8548859591470491
$ws.Range("M2").Value = 0.06720333333333334
$ws.Range("N2").Value = 0.20161
$ws.Range("Q2").Value = 0.1703976134433333
$ws.Range("R2").Value = 1.53357852099
$ws.Range("S2").Value = 0.8548859591470492
$ws.Range("T2").Value = 0.8548859591470491

# Update TPM-derived values for the MuSCs -> Thpo/Mpl -> ECs edge (now row 3)
$ws.Range("I3").Value = 0.1451140408529508
$ws.Range("J3").Value = 0.1451140408529508
$ws.Range("M3").Value = 0.06720333333333334
$ws.Range("N3").Value = 0.20161
$ws.Range("Q3").Value = 0.02892442667222222
$ws.Range("R3").Value = 0.26031984005
$ws.Range("S3").Value = 0.1451140408529508
$ws.Range("T3").Value = 0.1451140408529508
